# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted for "Apio" (Vega Monumental
# Concepción) at row 113, pushing the existing rows 113-170 down to 114-171.
# The new row carries its own date/price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 113 - this shifts rows 113:170 down to 114:171
# and Excel copies the formatting (incl. the date style on column D) from the
# row above, matching the target workbook.
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new weekly record.
$ws.Cells.Item(113, 1).Value = 11
$ws.Cells.Item(113, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(113, 3).Value = "Bíobío"
$ws.Cells.Item(113, 4).Value = 44523
$ws.Cells.Item(113, 5).Value = 8
$ws.Cells.Item(113, 6).Value = 100112017
$ws.Cells.Item(113, 7).Value = "Apio"
$ws.Cells.Item(113, 8).Value = "Americana (o)"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 100
$ws.Cells.Item(113, 11).Value = 6000
$ws.Cells.Item(113, 12).Value = 6500
$ws.Cells.Item(113, 13).Value = 6250
$ws.Cells.Item(113, 14).Value = "$/docena de matas"
$ws.Cells.Item(113, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(113, 16).Value = 1042
$ws.Cells.Item(113, 17).Value = 6
$ws.Cells.Item(113, 18).Value = "Hortaliza"
